# Commit: "modify data tables UI and modify important notes with last updated date"
#
# This adds a new "2024" column (column I) to the two region-level tables
# (region_soli = sheet index 3, region_asig = sheet index 4), and changes
# which sheet/view is currently active/selected in the workbook (UI state).

$wb = $excel.ActiveWorkbook

$wsMuniSoli = $wb.Worksheets.Item(1)   # municipios_soli
$wsMuniAsig = $wb.Worksheets.Item(2)   # municipios_asig
$wsRegSoli  = $wb.Worksheets.Item(3)   # region_soli
$wsRegAsig  = $wb.Worksheets.Item(4)   # region_asig

# ---------------------------------------------------------------------
# region_soli (sheet3): add 2024 data in column I
# ---------------------------------------------------------------------
$regSoli2024 = @(34, 38, 78, 27, 65, 39, 19, 19, 55, 69)

# Copy the formatting of the H1 header cell (bold + centered) onto I1
$wsRegSoli.Range("H1").Copy() | Out-Null
$wsRegSoli.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsRegSoli.Application.CutCopyMode = $false

$wsRegSoli.Range("I1").Value = 2024
for ($i = 0; $i -lt $regSoli2024.Length; $i++) {
    $row = $i + 2
    $wsRegSoli.Cells.Item($row, 9).Value = $regSoli2024[$i]
}

# ---------------------------------------------------------------------
# region_asig (sheet4): add 2024 data in column I
# ---------------------------------------------------------------------
$regAsig2024 = @(12, 16, 25, 2, 20, 17, 12, 11, 35, 38)

$wsRegAsig.Range("H1").Copy() | Out-Null
$wsRegAsig.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsRegAsig.Application.CutCopyMode = $false

$wsRegAsig.Range("I1").Value = 2024
for ($i = 0; $i -lt $regAsig2024.Length; $i++) {
    $row = $i + 2
    $wsRegAsig.Cells.Item($row, 9).Value = $regAsig2024[$i]
}

# ---------------------------------------------------------------------
# Update UI/selection state on each sheet, then leave region_soli as the
# active (front-most) sheet/tab, matching the saved workbook view.
# ---------------------------------------------------------------------
[void]$wsMuniSoli.Range("J2:J79").Select()
$wsMuniSoli.Activate() | Out-Null

[void]$wsMuniAsig.Range("L74").Select()
$wsMuniAsig.Activate() | Out-Null

[void]$wsRegAsig.Range("M13").Select()
$wsRegAsig.Activate() | Out-Null

[void]$wsRegSoli.Range("H21").Select()
$wsRegSoli.Activate() | Out-Null
